$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("R2").Value = 5124
$ws.Range("S2").Value = 458
$ws.Range("U2").Value = 32

# Row 3
$ws.Range("R3").Value = 903
$ws.Range("S3").Value = 61

# Row 4
$ws.Range("R4").Value = 260
$ws.Range("S4").Value = 24
$ws.Range("U4").Value = 3

# Row 6
$ws.Range("R6").Value = 297
$ws.Range("S6").Value = 39

# Row 7
$ws.Range("R7").Value = 936
$ws.Range("S7").Value = 115

# Row 11
$ws.Range("R11").Value = 1168
$ws.Range("S11").Value = 111
$ws.Range("U11").Value = 7

# Row 16 (Total row)
$ws.Range("R16").Value = 10895
$ws.Range("S16").Value = 942
$ws.Range("U16").Value = 73
